$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up the old content in rows 10-23 before rewriting the section ---
$ws.Range("A10:C23").ClearContents()

# --- Column layout cleanup: split the old merged col(1:2) range so col A stands alone ---
# (touching column B forces the engine to split the legacy "A:B" col-range without
#  disturbing column A's own stored width)
$origBWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(2).ColumnWidth = $origBWidth

# --- Row heights ---
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).EntireRow.AutoFit()
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Rows.Item(14).EntireRow.AutoFit()
$ws.Rows.Item(15).EntireRow.AutoFit()
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).RowHeight = 120
$ws.Rows.Item(20).EntireRow.AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 60
$ws.Rows.Item(24).RowHeight = 120
$ws.Rows.Item(25).EntireRow.AutoFit()
$ws.Rows.Item(26).RowHeight = 30

# --- Cell values ---
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "- Capacitar o aluno para relacionar as propriedades químicas e físicas dos elementos e seus compostos com suas posições na tabela periódica.-Capacitar o aluno a escrever os métodos industriais de obtenção dos elementos e seus compostos, bem como descrever suas aplicações- capacitar o aluno a comunicar-se eficazmente nas formas escrita, oral e gráfica"
$ws.Range("C10").Value = "- Capacitar o aluno para relacionar as propriedades químicas e físicas dos elementos e seus compostos com suas posições na tabela periódica.-Capacitar o aluno a escrever os métodos industriais de obtenção dos elementos e seus compostos, bem como descrever suas aplicações- capacitar o aluno a comunicar-se eficazmente nas formas escrita, oral e gráfica"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B13").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C13").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("B14").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C14").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("A16").Value = "Programa resumido:"
$ws.Range("B16").Value = "- Metais Representativos: Metais do Grupo 1, Metais do Grupo 2 e Metais do Grupo 13.- Metais de Transição: Propriedades gerais, Complexos."
$ws.Range("C16").Value = "- Metais Representativos: Metais do Grupo 1, Metais do Grupo 2 e Metais do Grupo 13.- Metais de Transição: Propriedades gerais, Complexos."
$ws.Range("A17").Value = "Short syllabus:"
$ws.Range("A18").Value = "Programa:"
$ws.Range("B18").Value = "Metais e compostos dos grupos 1, 2, 13 e de transição da Tabela Periódica: Propriedades físicas e químicas (relação com a posição na Tabela Periódica), processos de obtenção dos metais e compostos e aplicações - Formação de Complexos.Relacionar a disciplina com disciplinas anteriores e posteriores da grade do curso."
$ws.Range("C18").Value = "Metais e compostos dos grupos 1, 2, 13 e de transição da Tabela Periódica: Propriedades físicas e químicas (relação com a posição na Tabela Periódica), processos de obtenção dos metais e compostos e aplicações - Formação de Complexos.Relacionar a disciplina com disciplinas anteriores e posteriores da grade do curso."
$ws.Range("A19").Value = "Syllabus:"
$ws.Range("A20").Value = "Avaliação:"
$ws.Range("A21").Value = "Método:"
$ws.Range("B21").Value = "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."
$ws.Range("C21").Value = "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."
$ws.Range("A22").Value = "Critério:"
$ws.Range("B22").Value = "NF = (P1 + P2 + LE) /3"
$ws.Range("C22").Value = "NF = (P1 + P2 + LE) /3"
$ws.Range("A23").Value = "Norma de recuperação:"
$ws.Range("B23").Value = "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR"
$ws.Range("C23").Value = "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR"
$ws.Range("A24").Value = "Bibliografia:"
$ws.Range("B24").Value = "LEE, J. D. “Química Inorgânica não tão Concisa”, Editora Edgard Blücher, 1999. - SHRIVER, D. F.; ATKINS, P. W. “Química Inorgânica”, Editora Bookman, 4ª edição, 2008. - QUAGLIANO, J. V.; VALLARINO, L. “Química”, Editora Guanabara Koogan, 1973. - BUCHEL, K. H.; MORETTO, H. H.; WODITSCH, P. “Industrial Inorganic Chemistry”, Editora Wiley-VCH, 2000. - RAYNER-CANHAM, G.; OVERTON, T. “Química Inorgânica Descritiva”, Editora: Gen-LTC, 5ª edição, 2015. - SOUZA, M.M.V.M. “Processos Inorgânicos”, Editora: Synergia, 1ª edição, 2012."
$ws.Range("C24").Value = "LEE, J. D. “Química Inorgânica não tão Concisa”, Editora Edgard Blücher, 1999. - SHRIVER, D. F.; ATKINS, P. W. “Química Inorgânica”, Editora Bookman, 4ª edição, 2008. - QUAGLIANO, J. V.; VALLARINO, L. “Química”, Editora Guanabara Koogan, 1973. - BUCHEL, K. H.; MORETTO, H. H.; WODITSCH, P. “Industrial Inorganic Chemistry”, Editora Wiley-VCH, 2000. - RAYNER-CANHAM, G.; OVERTON, T. “Química Inorgânica Descritiva”, Editora: Gen-LTC, 5ª edição, 2015. - SOUZA, M.M.V.M. “Processos Inorgânicos”, Editora: Synergia, 1ª edição, 2012."
$ws.Range("A25").Value = "Requisitos:"
$ws.Range("B26").Value = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("C26").Value = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
